$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.179.33"
$ws.Range("E2").Value = "  +3.30%  "
$ws.Range("D3").Value = "3.076.02"
$ws.Range("E3").Value = "  +2.98%  "
$ws.Range("E4").Value = "  -0.05%  "
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "520.49"
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = "  +2.99%  "
$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.81"
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = "  +3.96%  "
$ws.Range("E7").Value = "  -0.04%  "
$origStyle = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.436"
$ws.Range("D8").Style = $origStyle
$ws.Range("E8").Value = "  +1.69%  "
$ws.Range("E9").Value = "  +2.01%  "
$ws.Range("E10").Value = "  +0.72%  "
$ws.Range("E11").Value = "  +3.32%  "
$ws.Range("D12").Value = "3.594.05"
$ws.Range("E12").Value = "  +2.79%  "
$ws.Range("E13").Value = "  +3.39%  "
$origStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.87"
$ws.Range("D14").Style = $origStyle
$ws.Range("E14").Value = "  +0.98%  "
$origStyle = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000164"
$ws.Range("D15").Style = $origStyle
$ws.Range("E15").Value = "  +1.16%  "
$ws.Range("D16").Value = "58.159.67"
$ws.Range("E16").Value = "  +3.36%  "
$ws.Range("D17").Value = "3.070.44"
$ws.Range("E17").Value = "  +2.51%  "
$origStyle = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.10"
$ws.Range("D18").Style = $origStyle
$ws.Range("E18").Value = "  +2.14%  "
$origStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.92"
$ws.Range("D19").Style = $origStyle
$ws.Range("E19").Value = "  +0.15%  "
$origStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.17"
$ws.Range("D20").Style = $origStyle
$ws.Range("E20").Value = "  +1.70%  "
$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "331.79"
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = "  +0.19%  "
$ws.Range("E22").Value = "  +0.13%  "
$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.500"
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = "  +1.69%  "
$origStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.06"
$ws.Range("D24").Style = $origStyle
$ws.Range("E24").Value = "  +2.28%  "
$ws.Range("E25").Value = "  +4.11%  "
$origStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = $origStyle
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").Value = "0.0₃0903"
$ws.Range("E27").Value = "  -2.60%  "
$origStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.45"
$ws.Range("D28").Style = $origStyle
$ws.Range("E28").Value = "  +1.72%  "
$origStyle = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.25"
$ws.Range("D29").Style = $origStyle
$ws.Range("E29").Value = "  +5.22%  "
$ws.Range("E30").Value = "  +2.72%  "
$origStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.20"
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = "  +3.72%  "
$origStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.74"
$ws.Range("D32").Style = $origStyle
$ws.Range("E32").Value = "  +2.65%  "
$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "155.28"
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = "  +1.84%  "
$origStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.53"
$ws.Range("D34").Style = $origStyle
$ws.Range("E34").Value = "  +1.47%  "
$origStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "27.20"
$ws.Range("D35").Style = $origStyle
$ws.Range("E35").Value = "  +4.57%  "
$ws.Range("E36").Value = "  +3.34%  "
$ws.Range("E37").Value = "  +1.48%  "
$origStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0676"
$ws.Range("D38").Style = $origStyle
$ws.Range("E38").Value = "  +2.66%  "
$ws.Range("D39").Value = "3.114.23"
$ws.Range("E39").Value = "  +3.01%  "
$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.94"
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = "  +4.14%  "
$origStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.75"
$ws.Range("D41").Style = $origStyle
$ws.Range("E41").Value = "  -0.38%  "
$origStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("D42").Style = $origStyle
$ws.Range("E42").Value = "  -0.04%  "
$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.654"
$ws.Range("D43").Style = $origStyle
$ws.Range("E43").Value = "  +0.45%  "
$ws.Range("D44").Value = "2.278.56"
$ws.Range("E44").Value = "  +4.41%  "
$ws.Range("E45").Value = "  +10.66%  "
$origStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "20.95"
$ws.Range("D46").Style = $origStyle
$ws.Range("E46").Value = "  +7.66%  "
$ws.Range("E47").Value = "  +1.73%  "
$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.90"
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = "  +1.37%  "
$origStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.931"
$ws.Range("D49").Style = $origStyle
$ws.Range("E49").Value = "  +1.25%  "
$origStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.732"
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = "  +9.16%  "
$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "261.18"
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = "  +14.05%  "
